$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect, make edits, then restore protection.
$ws.Unprotect()

# Update the confidential disclosure date (shared string in A42) from 2021-04-27 to 2021-04-28.
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

# Update the per-holding Weight (column D) and Percent Change (column E) figures for rows 2-39.

$ws.Range("D2").Value = 0.06262049289783393
$ws.Range("E2").Value = -0.006027234169208806
$ws.Range("D3").Value = 0.05677571706243909
$ws.Range("E3").Value = -0.02828568156659172
$ws.Range("D4").Value = 0.2962861503075179
$ws.Range("E4").Value = 0.003584229390680926
$ws.Range("D5").Value = 0.03703230117202183
$ws.Range("E5").Value = 0.01201780285185072
$ws.Range("D6").Value = 0.03238786243228801
$ws.Range("E6").Value = 0.01539732938976113
$ws.Range("D7").Value = 0.02950183030547208
$ws.Range("E7").Value = 0.006479338842975135
$ws.Range("D8").Value = 0.02829228235370099
$ws.Range("E8").Value = -0.007353842382645026
$ws.Range("D9").Value = 0.02399243799549665
$ws.Range("E9").Value = -0.003540974129209506
$ws.Range("D10").Value = 0.02482574956592486
$ws.Range("E10").Value = 0.02970781063125827
$ws.Range("D11").Value = 0.0230270319121488
$ws.Range("E11").Value = 0.01162829001548271
$ws.Range("D12").Value = 0.02247188819624455
$ws.Range("E12").Value = 0.003761283851554609
$ws.Range("D13").Value = 0.022008961404454
$ws.Range("E13").Value = -0.006769930675909919
$ws.Range("D14").Value = 0.02132259007973459
$ws.Range("E14").Value = 0.001677084921481997
$ws.Range("D15").Value = 0.02084405901406346
$ws.Range("E15").Value = -0.002744939018684445
$ws.Range("D16").Value = 0.02137460432600318
$ws.Range("E16").Value = 0.0002281368821290819
$ws.Range("D17").Value = 0.02108039874554641
$ws.Range("E17").Value = 0.01691212378235285
$ws.Range("D18").Value = 0.01507632928095397
$ws.Range("E18").Value = -0.006037605658099077
$ws.Range("D19").Value = 0.01647551250557934
$ws.Range("E19").Value = -0.003683241252302016
$ws.Range("D20").Value = 0.0152575122387896
$ws.Range("E20").Value = 0.001953125
$ws.Range("D21").Value = 0.01589316967339712
$ws.Range("E21").Value = 0.03013650062045747
$ws.Range("D22").Value = 0.01527354996472242
$ws.Range("E22").Value = -0.01467207764565659
$ws.Range("D23").Value = 0.0150958346233047
$ws.Range("E23").Value = 0.0001866368047780398
$ws.Range("D24").Value = 0.01433761862092678
$ws.Range("E24").Value = 0.006174845628859282
$ws.Range("D25").Value = 0.01391651995217723
$ws.Range("E25").Value = -0.005061319836480394
$ws.Range("D26").Value = 0.01450655655828667
$ws.Range("E26").Value = -0.002054231717337762
$ws.Range("D27").Value = 0.01288425388977163
$ws.Range("E27").Value = -0.01059722958140952
$ws.Range("D28").Value = 0.01335465772946328
$ws.Range("E28").Value = 0.02453748782862686
$ws.Range("D29").Value = 0.01416629669727958
$ws.Range("E29").Value = 0.009362808842652726
$ws.Range("D30").Value = 0.0130417053477472
$ws.Range("E30").Value = -0.005982451475671446
$ws.Range("D31").Value = 0.0123871927488673
$ws.Range("E31").Value = -0.007418293792427533
$ws.Range("D32").Value = 0.01328357159289619
$ws.Range("E32").Value = 0.00439698492462326
$ws.Range("D33").Value = 0.01257444403543427
$ws.Range("E33").Value = -0.003490175801447726
$ws.Range("D34").Value = 0.0066672511045171
$ws.Range("E34").Value = -0.006826271393046812
$ws.Range("D35").Value = 0.005478292125227331
$ws.Range("E35").Value = 0.001918702403323147
$ws.Range("D36").Value = 0.00582602903413554
$ws.Range("E36").Value = 0.008444312179153313
$ws.Range("D37").Value = 0.005605835391598464
$ws.Range("E37").Value = -0.003556792700842859
$ws.Range("D38").Value = 0.005053509114033758
$ws.Range("E38").Value = -0.01533183231478508
$ws.Range("D39").Value = 0.9999999999999998
$ws.Range("E39").Value = 0.001410452977983256

# Restore sheet protection (matches original legacy password hash "D382").
$ws.Protect("D382")

